$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 38171.4
$ws.Range("J3").Value = 38171.4
$ws.Range("L3").Value = 38171.4
$ws.Range("N3").Value = -38399.4
$ws.Range("H40").Value = 4100
$ws.Range("I40").Value = 3700
$ws.Range("K40").Value = 3700
$ws.Range("M40").Value = -3525
$ws.Range("H51").Value = 119659.78
$ws.Range("I51").Value = 258290.75
$ws.Range("J51").Value = 8755
$ws.Range("K51").Value = 258290.75
$ws.Range("L51").Value = 8755
$ws.Range("M51").Value = -257806.75
$ws.Range("N51").Value = -9723
$ws.Range("H95").Value = 18330
$ws.Range("J95").Value = 18330
$ws.Range("L95").Value = 18330
$ws.Range("N95").Value = -23822
$ws.Range("H102").Value = 38171.4
$ws.Range("J102").Value = 38171.4
$ws.Range("L102").Value = 38171.4
$ws.Range("N102").Value = -44661.4
$ws.Range("H106").Value = 8806.833000000001
$ws.Range("I106").Value = 2047.5454
$ws.Range("K106").Value = 2047.5454
$ws.Range("M106").Value = -1416.5454
$ws.Range("H111").Value = 1683.1111
$ws.Range("I111").Value = 1202.875
$ws.Range("K111").Value = 3608.625
$ws.Range("M111").Value = -541.625
$ws.Range("H121").Value = 2245
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 2490
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 7470
$ws.Range("M121").Value = -4253
$ws.Range("N121").Value = -10964
$ws.Range("H133").Value = 88000
$ws.Range("J133").Value = 88000
$ws.Range("L133").Value = 88000
$ws.Range("N133").Value = -98120
$ws.Range("H141").Value = 3885.7778
$ws.Range("I141").Value = 3885.7778
$ws.Range("K141").Value = 11657.3334
$ws.Range("M141").Value = -6477.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9436105
$ws.Range("I32").Value = 10418033
$ws.Range("K32").Value = 10418033
$ws.Range("M32").Value = -10417746
$ws.Range("H63").Value = 6474.143
$ws.Range("I63").Value = 4579.75
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 4579.75
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -3893.75
$ws.Range("N63").Value = -10372
$ws.Range("H66").Value = 6474.143
$ws.Range("I66").Value = 4579.75
$ws.Range("K66").Value = 22898.75
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -19466.75
$ws.Range("N66").Value = -51864
$ws.Range("H122").Value = 2354.0625
$ws.Range("I122").Value = 1066.5
$ws.Range("K122").Value = 3199.5
$ws.Range("M122").Value = -749.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3314
$ws.Range("I20").Value = 3409.8235
$ws.Range("J20").Value = 2499.5
$ws.Range("K20").Value = 3409.8235
$ws.Range("L20").Value = 2499.5
$ws.Range("M20").Value = -3162.8235
$ws.Range("N20").Value = -2993.5
$ws.Range("H86").Value = 3877.5715
$ws.Range("I86").Value = 3900
$ws.Range("K86").Value = 3900
$ws.Range("M86").Value = -2777
$ws.Range("H89").Value = 3877.5715
$ws.Range("I89").Value = 3900
$ws.Range("K89").Value = 19500
$ws.Range("M89").Value = -13884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 13205.125
$ws.Range("I5").Value = 316.5
$ws.Range("J5").Value = 26093.75
$ws.Range("K5").Value = 316.5
$ws.Range("L5").Value = 26093.75
$ws.Range("M5").Value = -204.5
$ws.Range("N5").Value = -26317.75
$ws.Range("H7").Value = 4854.7144
$ws.Range("I7").Value = 201.75
$ws.Range("J7").Value = 11058.667
$ws.Range("K7").Value = 201.75
$ws.Range("L7").Value = 11058.667
$ws.Range("M7").Value = -88.75
$ws.Range("N7").Value = -11284.667
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2826
$ws.Range("H22").Value = 442.16666
$ws.Range("I22").Value = 442.16666
$ws.Range("K22").Value = 442.16666
$ws.Range("M22").Value = -92.16665999999998
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 8958.571
$ws.Range("J55").Value = 7999.5
$ws.Range("L55").Value = 23998.5
$ws.Range("N55").Value = -24352.5
$ws.Range("H56").Value = 8830
$ws.Range("I56").Value = 8830
$ws.Range("K56").Value = 8830
$ws.Range("M56").Value = -8300
$ws.Range("H64").Value = 3799.6667
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12540
$ws.Range("H67").Value = 3799.6667
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13872
$ws.Range("H92").Value = 2505000
$ws.Range("H109").Value = 1086.6316
$ws.Range("I109").Value = 1086.6316
$ws.Range("K109").Value = 3259.8948
$ws.Range("M109").Value = -2219.8948
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 41787268
$ws.Range("I129").Value = 1207.25
$ws.Range("J129").Value = 83573330
$ws.Range("K129").Value = 3621.75
$ws.Range("L129").Value = 250719990
$ws.Range("M129").Value = 1378.25
$ws.Range("N129").Value = -250729990
$ws.Range("H134").Value = 8447.5
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -9930
$ws.Range("H141").Value = 344054.88
$ws.Range("J141").Value = 11785.571
$ws.Range("L141").Value = 35356.713
$ws.Range("N141").Value = -45716.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4882.9165
$ws.Range("I70").Value = 4781.364
$ws.Range("K70").Value = 4781.364
$ws.Range("M70").Value = -4511.364
$ws.Range("H73").Value = 4882.9165
$ws.Range("I73").Value = 4781.364
$ws.Range("K73").Value = 4781.364
$ws.Range("M73").Value = -3845.364
$ws.Range("H88").Value = 82000
$ws.Range("J88").Value = 82000
$ws.Range("L88").Value = 82000
$ws.Range("N88").Value = -82902
$ws.Range("H91").Value = 82000
$ws.Range("J91").Value = 82000
$ws.Range("L91").Value = 82000
$ws.Range("N91").Value = -85120
$ws.Range("H113").Value = 3864.0908
$ws.Range("I113").Value = 2811
$ws.Range("K113").Value = 2811
$ws.Range("M113").Value = -641
$ws.Range("H126").Value = 4148.467
$ws.Range("I126").Value = 3618.2856
$ws.Range("K126").Value = 10854.8568
$ws.Range("M126").Value = -8384.856800000001
$ws.Range("H132").Value = 142859820
$ws.Range("I132").Value = 200002940
$ws.Range("K132").Value = 600008820
$ws.Range("M132").Value = -600006290
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6912
$ws.Range("I46").Value = 4430.8887
$ws.Range("J46").Value = 9145
$ws.Range("K46").Value = 4430.8887
$ws.Range("L46").Value = 9145
$ws.Range("M46").Value = -4242.8887
$ws.Range("N46").Value = -9521
$ws.Range("H64").Value = 48874
$ws.Range("J64").Value = 48874
$ws.Range("L64").Value = 48874
$ws.Range("N64").Value = -49324
$ws.Range("H67").Value = 48874
$ws.Range("J67").Value = 48874
$ws.Range("L67").Value = 48874
$ws.Range("N67").Value = -50434
$ws.Range("H68").Value = 1719
$ws.Range("J68").Value = 1500
$ws.Range("L68").Value = 1500
$ws.Range("N68").Value = -2998
$ws.Range("H71").Value = 1719
$ws.Range("J71").Value = 1500
$ws.Range("L71").Value = 7500
$ws.Range("N71").Value = -14988
$ws.Range("H82").Value = 990.3333
$ws.Range("I82").Value = 612.6667
$ws.Range("J82").Value = 1745.6666
$ws.Range("K82").Value = 612.6667
$ws.Range("L82").Value = 1745.6666
$ws.Range("M82").Value = -251.6667
$ws.Range("N82").Value = -2467.6666
$ws.Range("H85").Value = 990.3333
$ws.Range("I85").Value = 612.6667
$ws.Range("J85").Value = 1745.6666
$ws.Range("K85").Value = 612.6667
$ws.Range("L85").Value = 1745.6666
$ws.Range("M85").Value = 635.3333
$ws.Range("N85").Value = -4241.6666
$ws.Range("H93").Value = 66668916
$ws.Range("I93").Value = 71430880
$ws.Range("J93").Value = 1499
$ws.Range("K93").Value = 71430880
$ws.Range("L93").Value = 1499
$ws.Range("M93").Value = -71429632
$ws.Range("N93").Value = -3995
$ws.Range("H97").Value = 20344
$ws.Range("J97").Value = 20344
$ws.Range("L97").Value = 20344
$ws.Range("N97").Value = -22326

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 89644
$ws.Range("J119").Value = 89644
$ws.Range("L119").Value = 89644
$ws.Range("N119").Value = -99320
$ws.Range("H126").Value = 5067.516
$ws.Range("I126").Value = 4662.2593
$ws.Range("J126").Value = 7803
$ws.Range("K126").Value = 13986.7779
$ws.Range("L126").Value = 23409
$ws.Range("M126").Value = -11516.7779
$ws.Range("N126").Value = -28349
$ws.Range("H136").Value = 1331.8
$ws.Range("I136").Value = 1314.2222
$ws.Range("J136").Value = 1490
$ws.Range("K136").Value = 3942.6666
$ws.Range("L136").Value = 4470
$ws.Range("M136").Value = -1392.6666
$ws.Range("N136").Value = -9570
